$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5: new wicket counts / matches played ---
$ws.Range("B2").Value = 21
$ws.Range("D2").Value = 10

$ws.Range("B3").Value = 28
$ws.Range("D3").Value = 10

$ws.Range("B4").Value = 19
$ws.Range("D4").Value = 10

$ws.Range("B5").Value = 23
$ws.Range("D5").Value = 10

# --- Add new bowlers (rows 6-10) ---
# Column A (names) for rows 6-9 first, then column C (teams) for rows 6-9,
# then row 10 (name then team) - matches the order the strings were typed in.
$ws.Range("A6").Value = "WARNE"
$ws.Range("A7").Value = "BOULT"
$ws.Range("A8").Value = "SHAKIB"
$ws.Range("A9").Value = "STARC"

$ws.Range("C6").Value = "RR"
$ws.Range("C7").Value = "DC"
$ws.Range("C8").Value = "KKR"
$ws.Range("C9").Value = "CSK"

$ws.Range("A10").Value = "JOHNSON"
$ws.Range("C10").Value = "RCB"

# --- Fill in wickets / matches played for the new rows ---
$ws.Range("B6").Value = 27
$ws.Range("D6").Value = 10

$ws.Range("B7").Value = 20
$ws.Range("D7").Value = 10

$ws.Range("B8").Value = 21
$ws.Range("D8").Value = 10

$ws.Range("B9").Value = 26
$ws.Range("D9").Value = 10

$ws.Range("B10").Value = 18
$ws.Range("D10").Value = 10

# --- Extend the wickets/match formula down through row 10 ---
$ws.Range("E6").Formula = "=B6/D6"
$ws.Range("E7").Formula = "=B7/D7"
$ws.Range("E8").Formula = "=B8/D8"
$ws.Range("E9").Formula = "=B9/D9"
$ws.Range("E10").Formula = "=B10/D10"

# --- Update the view: zoom + active selection cell ---
$ws.Range("B8").Select()
$excel.ActiveWindow.Zoom = 207
